# Update column C ("Förändrad") for rows 2-101 from 46074 to 46075 (+1 day).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 101; $r++) {
    $cell = $ws.Cells.Item($r, 3)  # Column C
    if ($cell.Value2 -eq 46074) {
        $cell.Value2 = 46075
    }
}
